{"js": "// Apply the eight Chai-Tee tagline corrections.\nconst replacements = [\n  [\"Chai Tee: Das Gew\u00fcrz des Lebens\", \"Chai-Tee: Das Gew\u00fcrz des Lebens\"],\n  [\"Chai-Tee: Eine Welt des Geschmacks in einer Tasse\", \"Chai-Tee: Eine Welt der Aromen in einer Tasse\"],\n  [\"Chai Tee: Die perfekte Mischung aus Gesundheit und Vergn\u00fcgen\", \"Chai-Tee: Die perfekte Mischung aus Gesundheit und Genuss\"],\n  [\"Chai Tee: Mehr als nur Tee, eine Lebensart\", \"Chai-Tee: Mehr als nur Tee, eine Lebensart\"],\n  [\"Chai Tee: Ein Getr\u00e4nk aus allen Jahreszeiten und Gr\u00fcnden\", \"Chai-Tee: Ein Getr\u00e4nk f\u00fcr alle Jahreszeiten und Anl\u00e4sse\"],\n  [\"Chai-Tee: Eine s\u00fc\u00dfe Flucht vom Alltag\", \"Chai-Tee: Eine s\u00fc\u00dfe Flucht aus dem Alltag\"],\n  [\"Chai-Tee: Teilen Sie die W\u00e4rme, teilen Sie die Liebe\", \"Chai-Tee: Gemeinsame W\u00e4rme, gemeinsame Liebe\"],\n  [\"Chai Tee: G\u00f6nnen Sie sich etwas Besonderes\", \"Chai-Tee: G\u00f6nnen Sie sich etwas Besonderes\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the eight Chai-Tee tagline corrections.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"Chai Tee: Das Gew\u00fcrz des Lebens\", \"Chai-Tee: Das Gew\u00fcrz des Lebens\"),\n  @(\"Chai-Tee: Eine Welt des Geschmacks in einer Tasse\", \"Chai-Tee: Eine Welt der Aromen in einer Tasse\"),\n  @(\"Chai Tee: Die perfekte Mischung aus Gesundheit und Vergn\u00fcgen\", \"Chai-Tee: Die perfekte Mischung aus Gesundheit und Genuss\"),\n  @(\"Chai Tee: Mehr als nur Tee, eine Lebensart\", \"Chai-Tee: Mehr als nur Tee, eine Lebensart\"),\n  @(\"Chai Tee: Ein Getr\u00e4nk aus allen Jahreszeiten und Gr\u00fcnden\", \"Chai-Tee: Ein Getr\u00e4nk f\u00fcr alle Jahreszeiten und Anl\u00e4sse\"),\n  @(\"Chai-Tee: Eine s\u00fc\u00dfe Flucht vom Alltag\", \"Chai-Tee: Eine s\u00fc\u00dfe Flucht aus dem Alltag\"),\n  @(\"Chai-Tee: Teilen Sie die W\u00e4rme, teilen Sie die Liebe\", \"Chai-Tee: Gemeinsame W\u00e4rme, gemeinsame Liebe\"),\n  @(\"Chai Tee: G\u00f6nnen Sie sich etwas Besonderes\", \"Chai-Tee: G\u00f6nnen Sie sich etwas Besonderes\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
